$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.936.73'
$ws.Range("E2").Value = '  -0.25%  '

$ws.Range("D3").Value = '2.587.58'
$ws.Range("E3").Value = '  +1.23%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.46'
$ws.Range("E5").Value = '  +0.20%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.12'
$ws.Range("E6").Value = '  -0.18%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.05%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.596'
$ws.Range("E8").Value = '  +2.08%  '

$ws.Range("E9").Value = '  +2.28%  '

$ws.Range("E10").Value = '  +2.32%  '

$ws.Range("E11").Value = '  -0.14%  '

$ws.Range("E12").Value = '  -0.23%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '27.38'
$ws.Range("E13").Value = '  +0.39%  '

$ws.Range("D14").Value = '3.050.55'
$ws.Range("E14").Value = '  +1.24%  '

$ws.Range("D15").Value = '62.829.37'
$ws.Range("E15").Value = '  -0.26%  '

$ws.Range("E16").Value = '  +2.92%  '

$ws.Range("D17").Value = '2.618.50'
$ws.Range("E17").Value = '  +2.59%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.31'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '342.24'
$ws.Range("E19").Value = '  +1.58%  '

$ws.Range("E20").Value = '  +0.76%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.70'
$ws.Range("E21").Value = '  -1.03%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  -0.10%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.06'
$ws.Range("E23").Value = '  +2.12%  '

$ws.Range("D24").Value = '2.709.11'
$ws.Range("E24").Value = '  +1.16%  '

$ws.Range("E25").Value = '  -1.65%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.60'
$ws.Range("E26").Value = '  -2.03%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.997'
$ws.Range("E27").Value = '  -0.29%  '

$ws.Range("B28").Value = 'Aptos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.84'
$ws.Range("E28").Value = '  +6.44%  '

$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.33'
$ws.Range("E29").Value = '  -0.65%  '

$ws.Range("E30").Value = '  -1.95%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.93'
$ws.Range("E31").Value = '  +0.06%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '470.86'
$ws.Range("E32").Value = '  +14.09%  '

$ws.Range("D33").Value = '0.0₃0823'
$ws.Range("E33").Value = '  +0.90%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '175.93'
$ws.Range("E34").Value = '  -0.73%  '

$ws.Range("E35").Value = '  +3.66%  '

$ws.Range("E36").Value = '  +0.05%  '

$ws.Range("E37").Value = '  +0.54%  '

$ws.Range("E38").Value = '  -0.61%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.53'
$ws.Range("E39").Value = '  +4.09%  '

$ws.Range("E41").Value = '  -2.04%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '157.63'
$ws.Range("E42").Value = '  +4.12%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.76'
$ws.Range("E43").Value = '  -0.16%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.635'
$ws.Range("E44").Value = '  +5.10%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.18'
$ws.Range("E45").Value = '  +1.04%  '

$ws.Range("E46").Value = '  +0.02%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0967'
$ws.Range("E47").Value = '  -0.24%  '

$ws.Range("E48").Value = '  -1.33%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '18.36'
$ws.Range("E49").Value = '  -0.19%  '

$ws.Range("E50").Value = '  -0.12%  '

$ws.Range("E51").Value = '  +0.88%  '
